$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Rectangle 5")
$tf = $shp.TextFrame
$tr = $tf.TextRange

# ---------------------------------------------------------------------------
# 1) Rewrite "portador do RG nº {{RG}} e CPF nº {{CPF}}" -> "portador do CPF nº {{CPF}}"
#    i.e. drop the RG mention entirely, keeping only the CPF mention.
#    Runs are edited surgically (Characters/Font tricks) so existing run
#    formatting (color/typefaces/size) is preserved and only the minimal
#    run split needed for "do " is introduced.
# ---------------------------------------------------------------------------
$full = $tr.Text
$blockStart = $full.IndexOf("portador do RG") + 1

# Split "do RG nº " away from "portador " (clean run boundary, no visual change)
$doStart = $blockStart + "portador ".Length
$doLen = "do RG nº ".Length
$doSplit = $tr.Characters($doStart, $doLen)
$doSplit.Font.Size = $doSplit.Font.Size

# Delete "RG nº " (tail of that split run), leaving "do " as its own run
$rgStart = $doStart + "do ".Length
$rgLen = "RG nº ".Length
$tr.Characters($rgStart, $rgLen).Text = ""

# Delete the "{{RG}}" run entirely (now located right after "do ")
$rgRunStart = $rgStart
$rgRunLen = "{{RG}}".Length
$tr.Characters($rgRunStart, $rgRunLen).Text = ""

# Replace the single-space run (formerly between {{RG}} and "e CPF nº ") with "CPF "
$spaceStart = $rgRunStart
$tr.Characters($spaceStart, 1).Text = "CPF "

# Delete the "e CPF " prefix of the following run, leaving "nº " as its own run
$afterCpfStart = $spaceStart + "CPF ".Length
$ecpfLen = "e CPF ".Length
$tr.Characters($afterCpfStart, $ecpfLen).Text = ""

# ---------------------------------------------------------------------------
# 2) Shrink the text box now that a line of text was removed.
#    (Only the height changes - position/width are left untouched.)
# ---------------------------------------------------------------------------
$shp.Height = 2523768 / 12700.0
